# Auto-generated edit script: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-leve sheets,
# matching a scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H17" = 2125.8914
    "J17" = 2125.8914
    "L17" = 6377.674199999999
    "N17" = -6713.674199999999
    "H62" = 12506166
    "I62" = 31253414
    "K62" = 31253414
    "M62" = -31252790
    "H65" = 12506166
    "I65" = 31253414
    "K65" = 156267070
    "M65" = -156263950
    "H116" = 10852.056
    "I116" = 5712
    "J116" = 15992.111
    "K116" = 5712
    "L116" = 15992.111
    "M116" = -2270
    "N116" = -22876.111
    "H133" = 38786.41
    "J133" = 38786.41
    "L133" = 38786.41
    "N133" = -48906.41
    "H137" = 3535.3845
    "I137" = 3061.375
    "J137" = 4293.8
    "K137" = 9184.125
    "L137" = 12881.4
    "M137" = -6634.125
    "N137" = -17981.4
    "H138" = 5508.244
    "J138" = 6213.485
    "L138" = 18640.455
    "N138" = -28920.455
    "H141" = 3497.8
    "I141" = 3497.8
    "J141" = 0
    "K141" = 10493.4
    "L141" = 0
    "M141" = -5313.400000000001
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H2" = 655.5
    "I2" = 655.3889
    "K2" = 655.3889
    "M2" = -542.3889
    "H32" = 1652.03
    "I32" = 1623.2626
    "K32" = 1623.2626
    "M32" = -1336.2626
    "H97" = 751
    "I97" = 798.8148
    "J97" = 105.5
    "K97" = 798.8148
    "L97" = 105.5
    "M97" = -302.8148
    "N97" = -1097.5
    "H116" = 655.5
    "I116" = 655.3889
    "K116" = 655.3889
    "M116" = 1638.6111
    "H122" = 3857.4517
    "J122" = 4395
    "L122" = 13185
    "N122" = -18085
    "H140" = 16750
    "J140" = 16750
    "L140" = 16750
    "N140" = -27110
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H3" = 655.5
    "I3" = 655.3889
    "K3" = 655.3889
    "M3" = -541.3889
    "H134" = 28914.684
    "I134" = 4125.1113
    "K134" = 12375.3339
    "M134" = -9840.333899999998
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 41798.93
    "I31" = 3682.625
    "K31" = 3682.625
    "M31" = -3387.625
    "H34" = 41798.93
    "I34" = 3682.625
    "K34" = 3682.625
    "M34" = -3480.625
    "H58" = 8015.5713
    "I58" = 5221.8
    "K58" = 5221.8
    "M58" = -5018.8
    "H59" = 30434.3
    "I59" = 0
    "J59" = 30434.3
    "K59" = 0
    "L59" = 30434.3
    "N59" = -32724.3
    "H86" = 5974.5
    "I86" = 4296.25
    "J86" = 7652.75
    "K86" = 4296.25
    "L86" = 7652.75
    "M86" = -3173.25
    "N86" = -9898.75
    "H89" = 5974.5
    "I89" = 4296.25
    "J89" = 7652.75
    "K89" = 21481.25
    "L89" = 38263.75
    "M89" = -15865.25
    "N89" = -49495.75
    "H99" = 4165.4443
    "I99" = 3898.2
    "K99" = 3898.2
    "M99" = -2400.2
    "H122" = 3873.4285
    "I122" = 3524
    "J122" = 4339.3335
    "K122" = 10572
    "L122" = 13018.0005
    "M122" = -8122
    "N122" = -17918.0005
    "H126" = 4165.4443
    "I126" = 3898.2
    "K126" = 11694.6
    "M126" = -9224.599999999999
    "H132" = 2296.3635
    "J132" = 2827
    "L132" = 8481
    "N132" = -13541
    "H136" = 8015.5713
    "I136" = 5221.8
    "K136" = 15665.4
    "M136" = -13115.4
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H36" = 140653.67
    "I36" = 946
    "J36" = 420069
    "K36" = 2838
    "L36" = 1260207
    "M36" = -2669
    "N36" = -1260545
    "H92" = 1000938.1
    "I92" = 1667147
    "J92" = 1624.75
    "K92" = 5001441
    "L92" = 4874.25
    "M92" = -5000193
    "N92" = -7370.25
    "H101" = 9536.556
    "I101" = 3300
    "J101" = 10316.125
    "K101" = 9900
    "L101" = 30948.375
    "M101" = -7466
    "N101" = -35816.375
    "H140" = 3991.875
    "I140" = 3775.8
    "J140" = 4352
    "K140" = 11327.4
    "L140" = 13056
    "M140" = -6147.400000000001
    "N140" = -23416
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H80" = 1003706.3
    "J80" = 1113854
    "L80" = 1113854
    "N80" = -1115850
    "H83" = 1003706.3
    "J83" = 1113854
    "L83" = 5569270
    "N83" = -5579254
    "H102" = 1676.3823
    "I102" = 1110.6428
    "K102" = 1110.6428
    "M102" = 511.3571999999999
    "H106" = 0
    "J106" = 0
    "L106" = 0
    "H122" = 4420.9287
    "I122" = 1713.2858
    "K122" = 5139.857400000001
    "M122" = -2689.857400000001
    "H138" = 52000
    "J138" = 52000
    "L138" = 52000
    "N138" = -62280
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H40" = 2977.8438
    "I40" = 3020.0967
    "K40" = 3020.0967
    "M40" = -2884.0967
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H104" = 25000
    "J104" = 25000
    "L104" = 25000
    "N104" = -31988
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

